$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("E2:E45")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq "fullRNASEQ") {
        $cell.Value = "fullRNASeq"
    }
}
